$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets in the right order / with the right
#    sheetId assignment:
#      - "Sheet1" is created FIRST (so it gets the lower internal sheetId)
#        but is positioned LAST (after VLAN_POOL).
#      - "VLAN_POOL" is created SECOND (higher sheetId) but inserted BEFORE
#        Sheet1, so visually it appears right after EPG.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet1 = $wb.Worksheets.Add($null, $lastSheet)
$sheet1.Name = "Sheet1"

$vlanPool = $wb.Worksheets.Add($sheet1, $null)
$vlanPool.Name = "VLAN_POOL"

# ---------------------------------------------------------------------------
# 2. Populate VLAN_POOL.
#    Cells are written in a specific order so that new shared-string
#    entries are introduced in the same sequence as the target workbook:
#    vlan_pool, vlan_pool_name, range_from, range_to, mark_pool,
#    allocation_mode.
# ---------------------------------------------------------------------------
$vlanPool.Range("A2").Value = "vlan_pool"
$vlanPool.Range("B1").Value = "vlan_pool_name"
$vlanPool.Range("C1").Value = "range_from"
$vlanPool.Range("D1").Value = "range_to"
$vlanPool.Range("B2").Value = "mark_pool"
$vlanPool.Range("E1").Value = "allocation_mode"
$vlanPool.Range("A1").Value = "type"
$vlanPool.Range("C2").Value = 1024
$vlanPool.Range("D2").Value = 1034

# Header row uses the existing bold header style.
$vlanPool.Range("A1:E1").Font.Bold = $true

# Body rows (2-22, columns A-G) get a left-aligned style.
$vlanPool.Range("A2:G22").HorizontalAlignment = -4131

# Column widths (character units; the host rounds to whole pixels using
# the workbook's Maximum Digit Width, so these inputs are chosen to land
# as close as possible to the authored widths).
$vlanPool.Columns.Item(1).ColumnWidth = 14.571428571428571
$vlanPool.Columns.Item(2).ColumnWidth = 21.142857142857142
$vlanPool.Columns.Item(3).ColumnWidth = 29.714285714285715
$vlanPool.Columns.Item(4).ColumnWidth = 19.428571428571427
$vlanPool.Columns.Item(5).ColumnWidth = 28.142857142857142

# ---------------------------------------------------------------------------
# 3. Sheet-view / selection tweaks on existing sheets.
# ---------------------------------------------------------------------------
# EPG: move the selection to B6 (requires briefly activating the sheet).
$epg = $wb.Worksheets.Item("EPG")
$epg.Activate()
$epg.Range("B6").Select()

# ---------------------------------------------------------------------------
# 4. Finally activate VLAN_POOL with E1 selected, and make it the sheet
#    that is active when the file is saved (this also clears the
#    tabSelected flag on whichever sheet was previously active, e.g. BD).
# ---------------------------------------------------------------------------
$vlanPool.Activate()
$vlanPool.Range("E1").Select()
